$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.90"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'22.74"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'6.212"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.06088"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'3.514"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'6.711"
$ws.Range("G7").Value = "'19"
$ws.Range("D8").Value = "'1.354"
$ws.Range("G8").Value = "'19"
$ws.Range("D9").Value = "'0.7991"
$ws.Range("G9").Value = "'19"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.08087"
$ws.Range("G11").Value = "'19"
$ws.Range("D12").Value = "'0.03339"
$ws.Range("G12").Value = "'19"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.09258"
$ws.Range("G14").Value = "'19"
$ws.Range("D15").Value = "'3.934"
$ws.Range("G15").Value = "'19"
$ws.Range("D16").Value = "'0.001717"
$ws.Range("G16").Value = "'19"
$ws.Range("D17").Value = "'0.04825"
$ws.Range("G17").Value = "'19"
$ws.Range("D18").Value = "'0.0006170"
$ws.Range("G18").Value = "'19"
$ws.Range("D19").Value = "'0.006201"
$ws.Range("G19").Value = "'19"
$ws.Range("D20").Value = "'0.001102"
$ws.Range("G20").Value = "'19"
$ws.Range("D21").Value = "'0.003388"
$ws.Range("G21").Value = "'19"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("G22").Value = "'19"
$ws.Range("D23").Value = "'3.696"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'2.296"
$ws.Range("G24").Value = "'19"
$ws.Range("D25").Value = "'0.3357"
$ws.Range("G25").Value = "'19"
$ws.Range("D26").Value = "'0.1253"
$ws.Range("G26").Value = "'19"
$ws.Range("D27").Value = "'0.0006166"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("G38").Value = "'19"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.04603"
$ws.Range("G40").Value = "'19"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007253"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'19"
$ws.Range("D42").Value = "'0.003902"
$ws.Range("G42").Value = "'19"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1120"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.01020"
$ws.Range("G44").Value = "'19"
$ws.Range("D45").Value = "'0.002970"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.00006008"
$ws.Range("G46").Value = "'19"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("G47").Value = "'19"
$ws.Range("D48").Value = "'0.7501"
$ws.Range("G48").Value = "'19"
$ws.Range("D49").Value = "'0.1297"
$ws.Range("G49").Value = "'19"
$ws.Range("D50").Value = "'0.00001501"
$ws.Range("G50").Value = "'19"
$ws.Range("D51").Value = "'0.01010"
$ws.Range("G51").Value = "'19"
